$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple single-cell "Taxonsorteringsordning" (column B) updates ---
$ws.Range("B2").Value = 57884
$ws.Range("B3").Value = 79243
$ws.Range("B4").Value = 79243
$ws.Range("B5").Value = 57884
$ws.Range("B6").Value = 57884
$ws.Range("B7").Value = 57884
$ws.Range("B8").Value = 57884
$ws.Range("B9").Value = 57884
$ws.Range("B14").Value = 79243
$ws.Range("B15").Value = 91828

# Date-like text (e.g. "2026-01-21") must stay plain text rather than be
# auto-converted into a date serial number by Excel, so force the number
# format of those date columns to Text before writing into them.
# (Each cell is set individually - a combined multi-area Range does not
# reliably apply NumberFormat to every area.)
$dateCells = "Y10","AA10","Y11","AA11","Y12","AA12","Y13","AA13"
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 10 becomes what row 11 used to contain (plus updated B id) ---
$ws.Range("A10").Value = 130894767
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("Q10").Value = 407194
$ws.Range("R10").Value = 7011100
$ws.Range("Y10").Value = "2026-01-21"
$ws.Range("Z10").Value = "12:26"
$ws.Range("AA10").Value = "2026-01-21"
$ws.Range("AB10").Value = "12:26"
$ws.Range("AC10").Value = "Ringhack"

# --- Row 11 becomes what row 10 used to contain (plus updated B id) ---
$ws.Range("A11").Value = 130894760
$ws.Range("B11").Value = 79243
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 406786
$ws.Range("R11").Value = 7010890
$ws.Range("Y11").Value = "2026-01-18"
$ws.Range("Z11").Value = "14:31"
$ws.Range("AA11").Value = "2026-01-18"
$ws.Range("AB11").Value = "14:31"
$ws.Range("AC11").Value = ""

# --- Row 12 swaps its observation identity with row 13 (plus updated B id) ---
$ws.Range("A12").Value = 130894766
$ws.Range("B12").Value = 57884
$ws.Range("Q12").Value = 407194
$ws.Range("R12").Value = 7011099
$ws.Range("Y12").Value = "2026-01-21"
$ws.Range("Z12").Value = "12:22"
$ws.Range("AA12").Value = "2026-01-21"
$ws.Range("AB12").Value = "12:22"

# --- Row 13 swaps its observation identity with row 12 (plus updated B id) ---
$ws.Range("A13").Value = 130894782
$ws.Range("B13").Value = 57884
$ws.Range("Q13").Value = 407192
$ws.Range("R13").Value = 7011093
$ws.Range("Y13").Value = "2026-01-24"
$ws.Range("Z13").Value = "15:03"
$ws.Range("AA13").Value = "2026-01-24"
$ws.Range("AB13").Value = "15:03"
